$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking strings so they remain shared strings
$ws.Range("A2:C2").NumberFormat = "@"

# Update row 2 values
$ws.Range("A2").Value = "456"
$ws.Range("B2").Value = "Teste"
$ws.Range("C2").Value = "154"

# Remove row 3 entirely (it's no longer present in the sheet)
$ws.Rows("3:3").Delete()
